$wb = $excel.ActiveWorkbook

# --- "Typography" sheet (sheet1) ---
$wsTypo = $wb.Worksheets.Item("Typography")

# Fnt_settings_digits (row 11): Wildcard Ranges column (I) changes from "0-9" to "a-z,A-Z,0-9"
$wsTypo.Range("I11").Value = "a-z,A-Z,0-9"

# Fnt_settings_total_digits (row 12): Size column (D) changes from 80 to 70
$wsTypo.Range("D12").Value = 70

# --- "Translation" sheet (sheet2) ---
$wsTrans = $wb.Worksheets.Item("Translation")

# Row 21 (SingleUseId30 / Fnt_Buttons): "SLOT #<value>" -> "SLOT#<value>"
$wsTrans.Range("F21").Value = "SLOT#<value>"

# Row 22 (SingleUseId31 / Fnt_Buttons): "1" -> "01"
# (leading apostrophe forces Excel to keep this as text instead of the number 1)
$wsTrans.Range("F22").Value = "'01"

# Row 26 (SingleUseId35 / Fnt_settings_digits): "0" -> "00"
$wsTrans.Range("F26").Value = "'00"

# Row 28 (SingleUseId37 / Fnt_settings_digits): "0" -> "00"
$wsTrans.Range("F28").Value = "'00"

# Row 29 (SingleUseId38 / Fnt_settings_digits): ":" -> "'"
# (doubled apostrophe: the first is consumed as Excel's text-prefix marker,
# leaving a literal single apostrophe character as the stored text)
$wsTrans.Range("F29").Value = "''"

# New rows 33-37
$wsTrans.Range("B33").Value = "SingleUseId42"
$wsTrans.Range("C33").Value = "Fnt_Buttons"
$wsTrans.Range("D33").Value = "Left"
$wsTrans.Range("E33").Value = "LTR"
$wsTrans.Range("F33").Value = "KEYZ"

$wsTrans.Range("B34").Value = "SingleUseId43"
$wsTrans.Range("C34").Value = "Fnt_Buttons"
$wsTrans.Range("D34").Value = "Left"
$wsTrans.Range("E34").Value = "LTR"
$wsTrans.Range("F34").Value = "NUM OF SLOTS:"

$wsTrans.Range("B35").Value = "SingleUseId44"
$wsTrans.Range("C35").Value = "Fnt_settings_digits"
$wsTrans.Range("D35").Value = "Left"
$wsTrans.Range("E35").Value = "LTR"
$wsTrans.Range("F35").Value = "<value>"

$wsTrans.Range("B36").Value = "SingleUseId45"
$wsTrans.Range("C36").Value = "Fnt_settings_digits"
$wsTrans.Range("D36").Value = "Left"
$wsTrans.Range("E36").Value = "LTR"
$wsTrans.Range("F36").Value = "IP DEV"

$wsTrans.Range("B37").Value = "SingleUseId46"
$wsTrans.Range("C37").Value = "Fnt_settings_digits"
$wsTrans.Range("D37").Value = "Left"
$wsTrans.Range("E37").Value = "LTR"
$wsTrans.Range("F37").Value = "`""
